$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-78)
# from serial date 45190 to 45192.
for ($r = 2; $r -le 78; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}
